$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("相談件数")

# The sheet currently ends with a single footnote cell on row 97 (column B).
# A new day's data row is being inserted as row 97, and the footnote moves
# down to row 98.

# 1) Push the footnote (row 97) down to row 98, preserving its content/format.
$ws.Rows.Item(97).Copy()
$ws.Rows.Item(98).PasteSpecial(-4104)   # xlPasteAll
$excel.CutCopyMode = $false

# 2) Clear the old footnote cell out of row 97 so it can hold the new data.
$ws.Range("A97:E97").ClearContents()

# 3) Copy the formatting of the previous data row (96) down onto row 97.
$ws.Range("A96:E96").Copy()
$ws.Range("A97:E97").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# 4) Fill in the new day's values (2020-05-01) on row 97.
$ws.Range("A97").Value = 43952
$ws.Range("B97").Value = 407
$ws.Range("C97").Value = 32436
$ws.Range("D97").Value = 109
$ws.Range("E97").Value = 6958

# 5) Refresh the active selection to the new bottom-right-most cell
#    (the frozen header pane - 1 row / 1 column - is kept as-is).
$ws.Activate()
[void]$ws.Range("E98").Select()
